# Apply "GP model updated 2" edits to Grain_lookup workbook.
# Sheets: Sheet1 (hidden), Wet_milling, Dry_grind

$wb = $excel.ActiveWorkbook
$wm = $wb.Worksheets.Item("Wet_milling")
$dg = $wb.Worksheets.Item("Dry_grind")

# ---------------------------------------------------------------------------
# Wet_milling sheet updates
# ---------------------------------------------------------------------------

# Row 12: P content in CGF without P recovery (mg/g) -> flat value of 12 for all scenarios
$wm.Range("D12:H12").Value = 12

# Row 26: previously blank "NA" row now becomes
# "P content in DDGS after P recovery (mg/g)" with a flat value of 2.5,
# formatted like the matching row on the Dry_grind sheet (numFmt 0.00, bold
# font, centered, bordered, shaded fill).
$dg.Range("D26:H26").Copy()
$wm.Range("D26:H26").PasteSpecial(-4122)   # xlPasteFormats
$wm.Range("D26:H26").Value = 2.5
$wm.Range("C26").Value = "P content in DDGS after P recovery (mg/g)"

# Row 37: Operating cost change for baseline (Delta $/yr) -> updated figures
$wm.Range("D37").Value = 1225061
$wm.Range("E37").Value = 1633141
$wm.Range("F37").Value = 2665990
$wm.Range("G37").Value = 3206420
$wm.Range("H37").Value = 6898226

# Row 40: Facility-dependent cost change (Delta $/yr) -> updated figures
$wm.Range("D40").Value = 478412
$wm.Range("E40").Value = 589762
$wm.Range("F40").Value = 826086
$wm.Range("G40").Value = 948702
$wm.Range("H40").Value = 1568332

# ---------------------------------------------------------------------------
# Dry_grind sheet updates
# ---------------------------------------------------------------------------

# Row 12: P content in DDGS without P recovery (mg/g) -> flat value of 9
$dg.Range("D12:H12").Value = 9

# Row 26: P content in DDGS after P recovery (mg/g) -> flat value of 3.2
$dg.Range("D26:H26").Value = 3.2

# Row 37: Operating cost change for baseline (Delta $/yr) -> updated figures
$dg.Range("D37").Value = 419192
$dg.Range("E37").Value = 738033
$dg.Range("F37").Value = 1058118
$dg.Range("G37").Value = 1818644
$dg.Range("H37").Value = 2134853

# Row 40: Facility-dependent cost change (Delta $/yr) -> updated figures
$dg.Range("D40").Value = 209076
$dg.Range("E40").Value = 324618
$dg.Range("F40").Value = 413174
$dg.Range("G40").Value = 604958
$dg.Range("H40").Value = 610722

# New column I (rows 2-41) mirroring the row-index notation column already
# used on the Wet_milling sheet (0-based sequential row index).
for ($r = 2; $r -le 41; $r++) {
    $dg.Cells.Item($r, 9).Value = $r - 2
}

# ---------------------------------------------------------------------------
# View / selection state: Wet_milling becomes the active sheet/tab, with a
# new selection; Dry_grind loses tabSelected and gets a new selection.
# ---------------------------------------------------------------------------
$dg.Range("I2:I41").Select()
$wm.Activate()
$wm.Range("F13").Select()
